$d = $word.ActiveDocument

$newText = "Ημερομηνίες παρατήρησης για τον αστερισμό του Leo: 14-23 Απριλίου, 14-23 Μαΐου"

# Collect the target paragraphs first (iterate over a snapshot of indices,
# since we will be mutating paragraph/run structure as we go).
$targets = @()
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Ημερομηνίες παρατήρησης*Περσεύς*") {
        $targets += $idx
    }
}

# Walk the collected paragraph indices from last to first so earlier
# indices/positions in the document stay valid while we edit.
for ($i = $targets.Count - 1; $i -ge 0; $i--) {
    $pIndex = $targets[$i]
    $p = $d.Paragraphs($pIndex)
    $start = $p.Range.Start
    $end = $p.Range.End - 1
    $r = $d.Range($start, $end)
    $r.Delete()
    $r2 = $d.Range($start, $start)
    $r2.InsertAfter($newText)
}
